$wb = $excel.ActiveWorkbook

# Grab the Analysis worksheet and copy its used range (A1:G65) as values
$wsAnalysis = $wb.Worksheets.Item("Analysis")
$wsAnalysis.Activate()
$wsAnalysis.Range("A1:G65").Copy()

# Add a new worksheet; it is added before the active sheet by default
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Sheet1"
$wsNew.Range("A1").PasteSpecial(-4104)

$wsAnalysis.Range("A1:G65").Copy()
$wsNew.Range("A1:G65").PasteSpecial(-4163)

# Move the new sheet so it sits right before "Terms only in one cluster"
$wsTerms = $wb.Worksheets.Item("Terms only in one cluster")
$wsNew.Move($wsTerms)

# Turn on AutoFilter for the pasted range on the new sheet
$wsNew.Range("A1:G65").AutoFilter()

$wsTerms.Activate()

$wb.Save()
